$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 16
$ws.Range("H16").Value = 2537.375
$ws.Range("I16").Value = 2383.1667
$ws.Range("K16").Value = 2383.1667
$ws.Range("M16").Value = -2153.1667
# Row 53
$ws.Range("H53").Value = 477.8
$ws.Range("I53").Value = 104.85714
$ws.Range("J53").Value = 804.125
$ws.Range("K53").Value = 104.85714
$ws.Range("L53").Value = 804.125
$ws.Range("M53").Value = 532.14286
$ws.Range("N53").Value = -2078.125
# Row 64
$ws.Range("H64").Value = 4941.1177
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 4999.9375
$ws.Range("K64").Value = 4000
$ws.Range("L64").Value = 4999.9375
$ws.Range("M64").Value = -3752
$ws.Range("N64").Value = -5495.9375
# Row 67
$ws.Range("H67").Value = 4941.1177
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 4999.9375
$ws.Range("K67").Value = 4000
$ws.Range("L67").Value = 4999.9375
$ws.Range("M67").Value = -3142
$ws.Range("N67").Value = -6715.9375
# Row 112
$ws.Range("H112").Value = 2216.9
$ws.Range("J112").Value = 2129.889
$ws.Range("L112").Value = 6389.667
$ws.Range("N112").Value = -8605.667000000001
# Row 125
$ws.Range("H125").Value = 2378.3572
$ws.Range("I125").Value = 2611.6667
$ws.Range("J125").Value = 2203.375
$ws.Range("K125").Value = 23505.0003
$ws.Range("L125").Value = 19830.375
$ws.Range("M125").Value = -21045.0003
$ws.Range("N125").Value = -24750.375
# Row 137
$ws.Range("H137").Value = 6494.4
$ws.Range("I137").Value = 2199.75
$ws.Range("K137").Value = 6599.25
$ws.Range("M137").Value = -4049.25

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1251.3636
$ws.Range("I2").Value = 1296.5
$ws.Range("K2").Value = 1296.5
$ws.Range("M2").Value = -1183.5
# Row 5
$ws.Range("H5").Value = 88.90000000000001
$ws.Range("J5").Value = 78.40000000000001
$ws.Range("L5").Value = 78.40000000000001
$ws.Range("N5").Value = -302.4
# Row 32
$ws.Range("H32").Value = 13164750
$ws.Range("I32").Value = 13520547
$ws.Range("J32").Value = 268
$ws.Range("K32").Value = 13520547
$ws.Range("L32").Value = 268
$ws.Range("M32").Value = -13520260
$ws.Range("N32").Value = -842
# Row 92
$ws.Range("H92").Value = 64179.668
$ws.Range("J92").Value = 64179.668
$ws.Range("L92").Value = 64179.668
$ws.Range("N92").Value = -69171.66800000001
# Row 102
$ws.Range("H102").Value = 2775.8
$ws.Range("I102").Value = 2510.6155
$ws.Range("K102").Value = 2510.6155
$ws.Range("M102").Value = -888.6154999999999
# Row 116
$ws.Range("H116").Value = 1251.3636
$ws.Range("I116").Value = 1296.5
$ws.Range("K116").Value = 1296.5
$ws.Range("M116").Value = 997.5
# Row 122
$ws.Range("H122").Value = 1420.7
$ws.Range("I122").Value = 950.3333
$ws.Range("K122").Value = 2850.9999
$ws.Range("M122").Value = -400.9998999999998
# Row 125
$ws.Range("H125").Value = 46340.8
$ws.Range("J125").Value = 46340.8
$ws.Range("L125").Value = 46340.8
$ws.Range("N125").Value = -56180.8
# Row 133
$ws.Range("H133").Value = 87999
$ws.Range("J133").Value = 87999
$ws.Range("L133").Value = 87999
$ws.Range("N133").Value = -93059

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1251.3636
$ws.Range("I3").Value = 1296.5
$ws.Range("K3").Value = 1296.5
$ws.Range("M3").Value = -1182.5
# Row 4
$ws.Range("H4").Value = 88.90000000000001
$ws.Range("J4").Value = 78.40000000000001
$ws.Range("L4").Value = 78.40000000000001
$ws.Range("N4").Value = -308.4
# Row 22
$ws.Range("H22").Value = 546
$ws.Range("I22").Value = 546
$ws.Range("K22").Value = 546
$ws.Range("M22").Value = -373
# Row 128
$ws.Range("H128").Value = 4688
$ws.Range("I128").Value = 4688
$ws.Range("K128").Value = 14064
$ws.Range("M128").Value = -11574

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 2636.7693
$ws.Range("I7").Value = 150.36363
$ws.Range("J7").Value = 16312
$ws.Range("K7").Value = 150.36363
$ws.Range("L7").Value = 16312
$ws.Range("M7").Value = -37.36363
$ws.Range("N7").Value = -16538
# Row 58
$ws.Range("H58").Value = 2021
$ws.Range("I58").Value = 1651.8
$ws.Range("K58").Value = 1651.8
$ws.Range("M58").Value = -1448.8
# Row 107
$ws.Range("H107").Value = 753.65216
$ws.Range("J107").Value = 1279.3334
$ws.Range("L107").Value = 1279.3334
$ws.Range("N107").Value = -5119.3334
# Row 136
$ws.Range("H136").Value = 2021
$ws.Range("I136").Value = 1651.8
$ws.Range("K136").Value = 4955.4
$ws.Range("M136").Value = -2405.4

$ws = $wb.Worksheets.Item("CUL")
# Row 86
$ws.Range("H86").Value = 96.25
$ws.Range("I86").Value = 95.333336
$ws.Range("J86").Value = 99
$ws.Range("K86").Value = 286.000008
$ws.Range("L86").Value = 297
$ws.Range("M86").Value = 899.999992
$ws.Range("N86").Value = -2669
# Row 89
$ws.Range("H89").Value = 96.25
$ws.Range("I89").Value = 95.333336
$ws.Range("J89").Value = 99
$ws.Range("K89").Value = 858.0000240000001
$ws.Range("L89").Value = 891
$ws.Range("M89").Value = 5069.999976
$ws.Range("N89").Value = -12747

$ws = $wb.Worksheets.Item("GSM")
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
# Row 70
$ws.Range("H70").Value = 4502
$ws.Range("I70").Value = 4502
$ws.Range("K70").Value = 4502
$ws.Range("M70").Value = -4232
# Row 73
$ws.Range("H73").Value = 4502
$ws.Range("I73").Value = 4502
$ws.Range("K73").Value = 4502
$ws.Range("M73").Value = -3566
# Row 102
$ws.Range("H102").Value = 2637.9
$ws.Range("I102").Value = 2653.2222
$ws.Range("K102").Value = 2653.2222
$ws.Range("M102").Value = -1031.2222
# Row 132
$ws.Range("H132").Value = 200002000
$ws.Range("I132").Value = 200002000
$ws.Range("K132").Value = 600006000
$ws.Range("M132").Value = -600003470

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 37228.17
$ws.Range("I7").Value = 2050.8333
$ws.Range("K7").Value = 2050.8333
$ws.Range("M7").Value = -1938.8333
# Row 24
$ws.Range("H24").Value = 80007
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
# Row 40
$ws.Range("H40").Value = 4036.8125
$ws.Range("I40").Value = 2499
$ws.Range("J40").Value = 5232.8887
$ws.Range("K40").Value = 2499
$ws.Range("L40").Value = 5232.8887
$ws.Range("M40").Value = -2363
$ws.Range("N40").Value = -5504.8887
# Row 126
$ws.Range("H126").Value = 37228.17
$ws.Range("I126").Value = 2050.8333
$ws.Range("K126").Value = 6152.499899999999
$ws.Range("M126").Value = -3682.499899999999
# Row 127
$ws.Range("H127").Value = 105920
$ws.Range("J127").Value = 105920
$ws.Range("L127").Value = 105920
$ws.Range("N127").Value = -115840
# Row 132
$ws.Range("H132").Value = 61881.473
$ws.Range("I132").Value = 32506.908
$ws.Range("K132").Value = 97520.724
$ws.Range("M132").Value = -94990.724
# Row 133
$ws.Range("H133").Value = 66666.664
$ws.Range("J133").Value = 66666.664
$ws.Range("L133").Value = 66666.664
$ws.Range("N133").Value = -71726.664

$ws = $wb.Worksheets.Item("WVR")
# Row 39
$ws.Range("H39").Value = 30495
$ws.Range("J39").Value = 30495
$ws.Range("L39").Value = 30495
$ws.Range("N39").Value = -31321
# Row 55
$ws.Range("H55").Value = 25149.4
$ws.Range("I55").Value = 31363
$ws.Range("J55").Value = 295
$ws.Range("K55").Value = 31363
$ws.Range("L55").Value = 295
$ws.Range("M55").Value = -31086
$ws.Range("N55").Value = -849
# Row 59
$ws.Range("H59").Value = 35000
$ws.Range("J59").Value = 35000
$ws.Range("L59").Value = 35000
$ws.Range("N59").Value = -36476
# Row 82
$ws.Range("H82").Value = 29666.334
$ws.Range("I82").Value = 29666.334
$ws.Range("K82").Value = 29666.334
$ws.Range("M82").Value = -29283.334
# Row 85
$ws.Range("H85").Value = 29666.334
$ws.Range("I85").Value = 29666.334
$ws.Range("K85").Value = 29666.334
$ws.Range("M85").Value = -28340.334
# Row 122
$ws.Range("H122").Value = 4270.032
$ws.Range("I122").Value = 3177.4736
$ws.Range("K122").Value = 9532.4208
$ws.Range("M122").Value = -7082.4208
# Row 126
$ws.Range("H126").Value = 1626.25
$ws.Range("I126").Value = 1608.5714
$ws.Range("K126").Value = 4825.7142
$ws.Range("M126").Value = -2355.7142
# Row 132
$ws.Range("H132").Value = 8810.429
$ws.Range("I132").Value = 1644.3
$ws.Range("J132").Value = 26725.75
$ws.Range("K132").Value = 4932.9
$ws.Range("L132").Value = 80177.25
$ws.Range("M132").Value = -2402.9
$ws.Range("N132").Value = -85237.25
